$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2739.1135
$ws.Range("J138").Value = 2995.0645
$ws.Range("L138").Value = 8985.193499999999
$ws.Range("N138").Value = -19265.1935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36639.305
$ws.Range("I32").Value = 37545.938
$ws.Range("K32").Value = 37545.938
$ws.Range("M32").Value = -37258.938

$ws.Range("H50").Value = 5783.3335
$ws.Range("I50").Value = 11800
$ws.Range("K50").Value = 11800
$ws.Range("M50").Value = -11086

$ws.Range("H74").Value = 1921.1578
$ws.Range("I74").Value = 1063.76
$ws.Range("K74").Value = 1063.76
$ws.Range("M74").Value = -189.76

$ws.Range("H77").Value = 1921.1578
$ws.Range("I77").Value = 1063.76
$ws.Range("K77").Value = 5318.8
$ws.Range("M77").Value = -950.8000000000002

$ws.Range("H122").Value = 2599.6365
$ws.Range("I122").Value = 1835.25
$ws.Range("K122").Value = 5505.75
$ws.Range("M122").Value = -3055.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 78331.664
$ws.Range("J58").Value = 78331.664
$ws.Range("L58").Value = 78331.664
$ws.Range("N58").Value = -78919.664

$ws.Range("H102").Value = 15150
$ws.Range("I102").Value = 15150
$ws.Range("K102").Value = 15150
$ws.Range("M102").Value = -11905

$ws.Range("H108").Value = 154015.67
$ws.Range("J108").Value = 154015.67
$ws.Range("L108").Value = 154015.67
$ws.Range("N108").Value = -161695.67

$ws.Range("H134").Value = 6064.0264
$ws.Range("I134").Value = 4588.4443
$ws.Range("J134").Value = 9685.909
$ws.Range("K134").Value = 13765.3329
$ws.Range("L134").Value = 29057.727
$ws.Range("M134").Value = -11230.3329
$ws.Range("N134").Value = -34127.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 881.625
$ws.Range("I22").Value = 297.5
$ws.Range("K22").Value = 297.5
$ws.Range("M22").Value = 52.5

$ws.Range("H96").Value = 24023.75
$ws.Range("J96").Value = 24023.75
$ws.Range("L96").Value = 24023.75
$ws.Range("N96").Value = -29515.75

$ws.Range("H117").Value = 87220.25
$ws.Range("J117").Value = 87220.25
$ws.Range("L117").Value = 87220.25
$ws.Range("N117").Value = -96398.25

$ws.Range("H134").Value = 4552.5293
$ws.Range("I134").Value = 2760.2856
$ws.Range("K134").Value = 8280.856800000001
$ws.Range("M134").Value = -5745.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 875.75
$ws.Range("J22").Value = 1251
$ws.Range("L22").Value = 3753
$ws.Range("N22").Value = -4091

$ws.Range("H27").Value = 875.75
$ws.Range("J27").Value = 1251
$ws.Range("L27").Value = 3753
$ws.Range("N27").Value = -3957

$ws.Range("H104").Value = 2900
$ws.Range("I104").Value = 2760
$ws.Range("J104").Value = 3250
$ws.Range("K104").Value = 8280
$ws.Range("L104").Value = 9750
$ws.Range("M104").Value = -5659
$ws.Range("N104").Value = -14992

$ws.Range("H114").Value = 713.8
$ws.Range("I114").Value = 731.44446
$ws.Range("J114").Value = 555
$ws.Range("K114").Value = 2194.33338
$ws.Range("L114").Value = 1665
$ws.Range("M114").Value = 1059.66662
$ws.Range("N114").Value = -8173

$ws.Range("H131").Value = 15158834
$ws.Range("J131").Value = 11036.571
$ws.Range("L131").Value = 33109.713
$ws.Range("N131").Value = -43189.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744

$ws.Range("H129").Value = 22000
$ws.Range("J129").Value = 22000
$ws.Range("L129").Value = 22000
$ws.Range("N129").Value = -32000

$ws.Range("H135").Value = 60999
$ws.Range("J135").Value = 60999
$ws.Range("L135").Value = 60999
$ws.Range("N135").Value = -71139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3810.2666
$ws.Range("I61").Value = 3456.25
$ws.Range("J61").Value = 4214.857
$ws.Range("K61").Value = 3456.25
$ws.Range("L61").Value = 4214.857
$ws.Range("M61").Value = -3254.25
$ws.Range("N61").Value = -4618.857

$ws.Range("H68").Value = 2465.2222
$ws.Range("I68").Value = 2465.2222
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2465.2222
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1716.2222
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 2465.2222
$ws.Range("I71").Value = 2465.2222
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 12326.111
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -8582.111000000001
$ws.Range("N71").Value = ""

$ws.Range("H113").Value = 3810.2666
$ws.Range("I113").Value = 3456.25
$ws.Range("J113").Value = 4214.857
$ws.Range("K113").Value = 3456.25
$ws.Range("L113").Value = 4214.857
$ws.Range("M113").Value = -1286.25
$ws.Range("N113").Value = -8554.857

$ws.Range("H122").Value = 5229.3125
$ws.Range("I122").Value = 5149.375
$ws.Range("K122").Value = 15448.125
$ws.Range("M122").Value = -12998.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 13496.333
$ws.Range("J61").Value = 20072.666
$ws.Range("L61").Value = 20072.666
$ws.Range("N61").Value = -20656.666

$ws.Range("H95").Value = 37988.668
$ws.Range("J95").Value = 37988.668
$ws.Range("L95").Value = 37988.668
$ws.Range("N95").Value = -43480.668

$ws.Range("H100").Value = 1290.579
$ws.Range("I100").Value = 1290.579
$ws.Range("K100").Value = 2581.158
$ws.Range("M100").Value = -2040.158

$ws.Range("H113").Value = 982.7895
$ws.Range("I113").Value = 1125
$ws.Range("K113").Value = 3375
$ws.Range("M113").Value = -1205

$ws.Range("H122").Value = 2899.5454
$ws.Range("I122").Value = 2811.4119
$ws.Range("K122").Value = 8434.235700000001
$ws.Range("M122").Value = -5984.235700000001

$ws.Range("H126").Value = 9408.333000000001
$ws.Range("I126").Value = 12481.667
$ws.Range("K126").Value = 37445.001
$ws.Range("M126").Value = -34975.001

$ws.Range("H136").Value = 3758.077
$ws.Range("I136").Value = 1618.0769
$ws.Range("J136").Value = 5898.077
$ws.Range("K136").Value = 4854.2307
$ws.Range("L136").Value = 17694.231
$ws.Range("M136").Value = -2304.2307
$ws.Range("N136").Value = -22794.231
